# "Fruta / hortaliza, semanal" update:
# A new weekly observation (row 117/118 data) is inserted at the top of the
# Betarraga price history block (rows 117-166), pushing every existing row
# down by one pair (2 rows), and the former last pair (old rows 165-166)
# becomes two brand-new rows (167-168) at the bottom.
#
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant per row (or depend only on
# odd/even row = Primera/Segunda), so only D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado) and
# P (Precio $/Kg) actually change value across the shift.
#
# Note: reading via the plain `.Value` property getter does not reliably
# round-trip through this COM host, so all reads below use `.Value2`.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, materialise the two brand-new rows (167, 168) by copying the
# non-shifting columns from what is currently the last pair (165, 166),
# since those rows don't exist yet in the sheet.
for ($col = 1; $col -le 18; $col++) {
    if ($col -eq 4 -or $col -eq 10 -or $col -eq 11 -or $col -eq 12 -or $col -eq 13 -or $col -eq 16) {
        continue
    }
    $v1 = $ws.Cells.Item(165, $col).Value2
    $ws.Cells.Item(167, $col).Value2 = $v1
    $v2 = $ws.Cells.Item(166, $col).Value2
    $ws.Cells.Item(168, $col).Value2 = $v2
}

# Shift the value-bearing columns (D,J,K,L,M,P) down by 2 rows, from the
# bottom up so a row's source data isn't clobbered before it's read.
for ($r = 168; $r -ge 119; $r--) {
    $src = $r - 2

    $dVal = $ws.Cells.Item($src, 4).Value2
    $dFmt = $ws.Cells.Item($src, 4).NumberFormat
    $ws.Cells.Item($r, 4).Value2 = $dVal
    $ws.Cells.Item($r, 4).NumberFormat = $dFmt

    $jVal = $ws.Cells.Item($src, 10).Value2
    $ws.Cells.Item($r, 10).Value2 = $jVal

    $kVal = $ws.Cells.Item($src, 11).Value2
    $ws.Cells.Item($r, 11).Value2 = $kVal

    $lVal = $ws.Cells.Item($src, 12).Value2
    $ws.Cells.Item($r, 12).Value2 = $lVal

    $mVal = $ws.Cells.Item($src, 13).Value2
    $ws.Cells.Item($r, 13).Value2 = $mVal

    $pVal = $ws.Cells.Item($src, 16).Value2
    $ws.Cells.Item($r, 16).Value2 = $pVal
}

# Finally, write the brand-new observation into rows 117/118 (the values
# that pushed everything else down).
$dateFmt = $ws.Cells.Item(119, 4).NumberFormat

$ws.Cells.Item(117, 4).Value2 = 44466
$ws.Cells.Item(117, 4).NumberFormat = $dateFmt
$ws.Cells.Item(117, 10).Value2 = 1200
$ws.Cells.Item(117, 11).Value2 = 400
$ws.Cells.Item(117, 12).Value2 = 450
$ws.Cells.Item(117, 13).Value2 = 425
$ws.Cells.Item(117, 16).Value2 = 106

$ws.Cells.Item(118, 4).Value2 = 44466
$ws.Cells.Item(118, 4).NumberFormat = $dateFmt
$ws.Cells.Item(118, 10).Value2 = 1200
$ws.Cells.Item(118, 11).Value2 = 400
$ws.Cells.Item(118, 12).Value2 = 450
$ws.Cells.Item(118, 13).Value2 = 425
$ws.Cells.Item(118, 16).Value2 = 85
